# Add Ashley King (White team) to the roster, matching the formatting
# she already had when her info was pasted in (dark "Segoe UI Historic"
# run copied from a Messenger/Facebook conversation).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 83

# -- values -------------------------------------------------------------
$ws.Cells.Item($row, 1).Value = "White"
$ws.Cells.Item($row, 2).Value = "Ashley King"
$ws.Cells.Item($row, 3).Value = "King, Ashley"
$ws.Cells.Item($row, 4).Value = "Female"
$ws.Cells.Item($row, 5).Value = "Player"
$ws.Cells.Item($row, 6).Value = "aking6@highpoint.edu"
$ws.Cells.Item($row, 7).Value = "M"

# -- plain-column formatting (matches the rest of the table) ------------
# A82/C82/E82 already carry the workbook's common "plain" cell style;
# copy that format onto the new row's non-email columns so they end up
# sharing the exact same style index instead of minting new ones.
$plain = $ws.Range("A82")
$plain.Copy()
$ws.Range("A83:E83").PasteSpecial(-4122)
$ws.Range("G83").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# -- special formatting for the pasted-in email cell ---------------------
# This run still has the dark Messenger/Facebook paste formatting
# (Segoe UI Historic, light-gray text on a dark-gray fill) from the
# original clipboard paste.
$email = $ws.Cells.Item($row, 6)
$email.Font.Name = "Segoe UI Historic"
$email.Font.Size = 11
$email.Font.Color = 15460068
$email.Interior.Color = 4341822
$email.Interior.PatternColor = 4341822

Write-Host "Added Ashley King to row $row"
